$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric/percentage need to be forced to
# text (matching the workbook's inlineStr storage) so Excel does not
# auto-convert them to numbers when the value is assigned. Multi-area
# Range().NumberFormat only reliably applies to the first area here, so
# loop per cell instead.
$textCells = @(
    "E2",
    "E3",
    "E4",
    "D5",
    "E5",
    "E6",
    "E7",
    "D8",
    "E8",
    "E9",
    "E10",
    "E11",
    "E12",
    "E13",
    "E14",
    "E15",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "E20",
    "D22",
    "E22",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "E26",
    "E27",
    "D28",
    "E28",
    "D29",
    "E29",
    "E30",
    "E31",
    "E32",
    "E33",
    "E34",
    "E35",
    "E36",
    "D37",
    "E37",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "E41",
    "E42",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "E47",
    "D48",
    "E48",
    "E49",
    "D50",
    "E50",
    "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.270.26'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '1.575.46'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '207.92'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('E6').Value = '  -1.83%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '22.29'
$ws.Range('E8').Value = '  +0.25%  '
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('E11').Value = '  -0.13%  '
$ws.Range('D12').Value = '1.800.27'
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D13').Value = '1.577.63'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').Value = '27.287.84'
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').Value = '62.49'
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('D18').Value = '215.03'
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').Value = '7.36'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('E20').Value = '  -0.75%  '
$ws.Range('D22').Value = '4.13'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('D24').Value = '2.00'
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('D25').Value = '152.20'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('E26').Value = '  -3.97%  '
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '0.104'
$ws.Range('E28').Value = '  -0.84%  '
$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('E30').Value = '  -1.63%  '
$ws.Range('E31').Value = '  -1.83%  '
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('D33').Value = '1.411.50'
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('E34').Value = '  -1.37%  '
$ws.Range('E35').Value = '  +1.69%  '
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('D37').Value = '0.938'
$ws.Range('E37').Value = '  -2.75%  '
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('D39').Value = '0.820'
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('D40').Value = '0.520'
$ws.Range('E40').Value = '  -2.54%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('E42').Value = '  +2.45%  '
$ws.Range('E43').Value = '  +3.70%  '
$ws.Range('D44').Value = '5.35'
$ws.Range('E44').Value = '  +1.90%  '
$ws.Range('D45').Value = '63.93'
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('D47').Value = '1.712.20'
$ws.Range('E47').Value = '  -0.95%  '
$ws.Range('D48').Value = '86.14'
$ws.Range('E48').Value = '  +0.38%  '
$ws.Range('D49').Value = '0.0₇0991'
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('D50').Value = '0.0954'
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('E51').Value = '  -0.03%  '

# Restore the default (unstyled) cell style now that the text values are
# committed as text -- keeps styling identical to the original cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
